# Updates the cryptocurrency price/volume snapshot on Sheet1.
# - Refreshes Price (D) and Volume(1h) (E) figures for every listed coin.
# - A handful of coins (rows 44-50) changed rank order, so their
#   Coin name (B) and Link (C) are rewritten too.
#
# Price values keep their original "text, not number" representation
# (the source feed renders already-formatted strings like "0.9989" or
# "29.519.94" into inline strings). A leading apostrophe is used where a
# value would otherwise be auto-parsed as a literal number by Excel, so
# the stored cell type/format stays exactly as authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.519.94"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.841.86"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'244.59"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'0.6282"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.07446"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "'0.2948"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'23.69"
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("D11").Value = "'0.07670"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.840.10"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "'5.025"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'0.6787"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "'83.97"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "'0.000009358"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "'5.971"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "29.499.90"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "2.080.18"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'237.24"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "'12.59"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "'7.357"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "'159.25"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "'0.1415"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'8.528"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'17.80"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'0.06075"
$ws.Range("E29").Value = "  +9.95%  "
$ws.Range("D30").Value = "'1.497"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "'1.246"
$ws.Range("E31").Value = "  +3.01%  "
$ws.Range("D32").Value = "'4.124"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "'4.101"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "'1.876"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("D35").Value = "'1.145"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").Value = "'0.7289"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'2.615"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "1.221.21"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("D40").Value = "'0.01768"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "'6.292"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").Value = "'0.9161"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.002.10"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'102.04"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "'65.68"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("D47").Value = "'0.5080"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.278"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.4067"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("D51").Value = "'0.1141"
$ws.Range("E51").Value = "  +4.14%  "
